{"js": "// Replace the 25 answer cells (5 columns x 5 data rows, every 4th table row)\n// in document order with their new values, per the commit diff. The mapping\n// is strictly positional \u2014 several old values repeat (e.g. \"20\u00f78=2, 4\",\n// \"42\u00f79=4, 6\") but map to different new values depending on where they sit,\n// so we must address each cell by (row, column) rather than doing a global\n// find-and-replace.\nconst replacements = [\n  \"89\u00f74=22, 1\", \"79\u00f78=9, 7\", \"30\u00f79=3, 3\", \"65\u00f74=16, 1\", \"89\u00f73=29, 2\",\n  \"23\u00f79=2, 5\", \"93\u00f77=13, 2\", \"15\u00f73=5, 0\", \"97\u00f76=16, 1\", \"13\u00f76=2, 1\",\n  \"61\u00f72=30, 1\", \"85\u00f75=17, 0\", \"86\u00f76=14, 2\", \"93\u00f77=13, 2\", \"55\u00f72=27, 1\",\n  \"45\u00f76=7, 3\", \"74\u00f77=10, 4\", \"38\u00f75=7, 3\", \"62\u00f75=12, 2\", \"49\u00f72=24, 1\",\n  \"80\u00f73=26, 2\", \"81\u00f73=27, 0\", \"46\u00f79=5, 1\", \"55\u00f78=6, 7\", \"46\u00f78=5, 6\",\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Data lives in every 4th row (0, 4, 8, 12, 16), 5 cells each = 25 cells.\nconst dataRowIndices = [];\nfor (let r = 0; r < rows.items.length; r += 4) {\n  dataRowIndices.push(r);\n}\n\nlet k = 0;\nfor (const r of dataRowIndices) {\n  const row = rows.items[r];\n  const cells = row.cells;\n  cells.load(\"items\");\n  await context.sync();\n  for (let c = 0; c < cells.items.length; c++) {\n    if (k >= replacements.length) break;\n    const cell = cells.items[c];\n    const paragraphs = cell.body.paragraphs;\n    paragraphs.load(\"items\");\n    await context.sync();\n    const para = paragraphs.items[0];\n    para.insertText(replacements[k], \"Replace\");\n    k++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 answer cells (5 columns x 5 data rows, every 4th table\n# row) in document order with their new values, per the commit diff. The\n# mapping is strictly positional -- several old values repeat (e.g.\n# \"20\u00f78=2, 4\", \"42\u00f79=4, 6\") but map to different new values depending on\n# where they sit, so each cell is addressed by its (row, column) position\n# rather than via a global find-and-replace.\n$replacements = @(\n  \"89\u00f74=22, 1\", \"79\u00f78=9, 7\", \"30\u00f79=3, 3\", \"65\u00f74=16, 1\", \"89\u00f73=29, 2\",\n  \"23\u00f79=2, 5\", \"93\u00f77=13, 2\", \"15\u00f73=5, 0\", \"97\u00f76=16, 1\", \"13\u00f76=2, 1\",\n  \"61\u00f72=30, 1\", \"85\u00f75=17, 0\", \"86\u00f76=14, 2\", \"93\u00f77=13, 2\", \"55\u00f72=27, 1\",\n  \"45\u00f76=7, 3\", \"74\u00f77=10, 4\", \"38\u00f75=7, 3\", \"62\u00f75=12, 2\", \"49\u00f72=24, 1\",\n  \"80\u00f73=26, 2\", \"81\u00f73=27, 0\", \"46\u00f79=5, 1\", \"55\u00f78=6, 7\", \"46\u00f78=5, 6\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$k = 0\n# Data lives in every 4th row (1-based rows 1, 5, 9, 13, 17), 5 columns each.\nfor ($r = 1; $r -le $t.Rows.Count; $r += 4) {\n  for ($c = 1; $c -le $t.Columns.Count; $c++) {\n    $t.Cell($r, $c).Range.Text = $replacements[$k]\n    $k++\n  }\n}\n"}
